$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows before the old row 1087 (2019-11-29), shifting the
# remaining data (old rows 1087-1159) down to rows 1096-1168. This fills
# in the missing trading days between 2019-11-15 and 2019-11-29.
$ws.Range("A1087:I1095").EntireRow.Insert()

# New row data: timestamp, date(text), id(text), name, open, high, low, close, vol
$newRows = @(
    @{ Row=1087; A=1574035200; B="2019-11-18"; C="5265"; D="DOLPHIN"; E=0.125; F=0.13;  G=0.125; H=0.13;  I=322000 },
    @{ Row=1088; A=1574121600; B="2019-11-19"; C="5265"; D="DOLPHIN"; E=0.125; F=0.125; G=0.12;  H=0.12;  I=564000 },
    @{ Row=1089; A=1574208000; B="2019-11-20"; C="5265"; D="DOLPHIN"; E=0.12;  F=0.12;  G=0.115; H=0.115; I=370600 },
    @{ Row=1090; A=1574294400; B="2019-11-21"; C="5265"; D="DOLPHIN"; E=0.115; F=0.125; G=0.115; H=0.12;  I=609500 },
    @{ Row=1091; A=1574380800; B="2019-11-22"; C="5265"; D="DOLPHIN"; E=0.125; F=0.13;  G=0.12;  H=0.12;  I=321200 },
    @{ Row=1092; A=1574640000; B="2019-11-25"; C="5265"; D="DOLPHIN"; E=0.12;  F=0.13;  G=0.12;  H=0.13;  I=735400 },
    @{ Row=1093; A=1574726400; B="2019-11-26"; C="5265"; D="DOLPHIN"; E=0.13;  F=0.135; G=0.125; H=0.125; I=735100 },
    @{ Row=1094; A=1574812800; B="2019-11-27"; C="5265"; D="DOLPHIN"; E=0.125; F=0.125; G=0.12;  H=0.125; I=343000 },
    @{ Row=1095; A=1574899200; B="2019-11-28"; C="5265"; D="DOLPHIN"; E=0.125; F=0.13;  G=0.12;  H=0.13;  I=391000 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = "'" + $r.B
    $ws.Cells.Item($row, 3).Value = "'" + $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
}
